# "added 4wk low sales check"
#
# Updates the per-week forecast/trend/inventory metrics on the
# "Forecast Comparison" sheet (rows 2-17, weeks W10-W25) to reflect a new
# "High Volume Season" trend classification with refreshed Inventory
# Coverage / Stockout Risk / Reorder Urgency / Seasonality Index numbers,
# and refreshes the dependent roll-up figures on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Forecast Comparison")

# Helper: write a text value into a cell while preserving the default
# (un-styled) cell format - without this, a numeric-looking string like
# "57" gets auto-coerced to a number by Excel's normal typed-input rules.
function Set-TextCell($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# r => @{ D=newForecast; G=newTrend; H=newCoverage; I=newStockoutRisk; J=newReorderUrgency; L=newSeasonality }
$updates = @{
    2  = @{ D = 102; G = "High Volume Season"; H = 10.8;  L = 1.19 }
    3  = @{ D = 206; G = "High Volume Season"; H = 4.86;  L = 1.09 }
    4  = @{ D = 350; G = "High Volume Season"; H = 2.27;  L = 1.11 }
    5  = @{ D = 430; G = "High Volume Season"; H = 1.04;  L = 1.1 }
    6  = @{ D = 344; G = "High Volume Season"; H = 0.04;  I = "High"; J = "Urgent"; L = 0.84 }
    7  = @{ D = 110; G = "High Volume Season"; H = 0;     I = "High"; J = "Urgent"; L = 1.12 }
    8  = @{          G = "High Volume Season"; H = 0;     I = "High"; J = "Urgent"; L = 1.04 }
    9  = @{ D = 57;  G = "High Volume Season"; H = 0;     I = "High"; J = "Urgent"; L = 1.16 }
    10 = @{ D = 79;  G = "High Volume Season"; H = 0;     I = "High"; J = "Urgent"; L = 0.96 }
    11 = @{ D = 334; G = "High Volume Season"; H = 0;     I = "High"; J = "Urgent"; L = 0.97 }
    12 = @{ D = 399; G = "High Volume Season"; H = 0;     I = "High"; J = "Urgent"; L = 0.85 }
    13 = @{ D = 207; G = "High Volume Season"; H = 0;     I = "High"; J = "Urgent"; L = 1.14 }
    14 = @{ D = 57;  G = "High Volume Season"; H = 0;     I = "High"; J = "Urgent"; L = 0.87 }
    15 = @{ D = 57;  G = "High Volume Season"; H = 0;     I = "High"; J = "Urgent"; L = 1.04 }
    16 = @{ D = 57;  G = "High Volume Season"; H = 0;     L = 1.16 }
    17 = @{ D = 322; G = "High Volume Season";            L = 0.82 }
}

foreach ($r in $updates.Keys) {
    $row = $updates[$r]

    if ($row.ContainsKey("D")) { $ws.Cells.Item($r, 4).Value = $row["D"] }
    if ($row.ContainsKey("G")) { $ws.Cells.Item($r, 7).Value = $row["G"] }
    if ($row.ContainsKey("H")) { $ws.Cells.Item($r, 8).Value = $row["H"] }
    if ($row.ContainsKey("I")) { $ws.Cells.Item($r, 9).Value = $row["I"] }
    if ($row.ContainsKey("J")) { $ws.Cells.Item($r, 10).Value = $row["J"] }
    if ($row.ContainsKey("L")) { $ws.Cells.Item($r, 12).Value = $row["L"] }
}

# Refresh the dependent roll-up numbers on the Summary sheet. These are
# stored as text cells ("1234" etc.), not numbers, so use Set-TextCell to
# keep their type/format unchanged.
$summary = $wb.Worksheets.Item("Summary")

Set-TextCell $summary.Cells.Item(9, 2)  "3176"   # Total Forecast (16 Weeks)
Set-TextCell $summary.Cells.Item(10, 2) "1659"   # Total Forecast (8 Weeks)
Set-TextCell $summary.Cells.Item(11, 2) "1090"   # Total Forecast (4 Weeks)
Set-TextCell $summary.Cells.Item(12, 2) "430"    # Max Forecast
Set-TextCell $summary.Cells.Item(14, 2) "58"     # Min Forecast
